# Update the cached "datetimeFigureOut" date placeholder text from
# 14/04/2025 to 18/04/2025 across the slide master, every slide layout,
# and the notes master (the placeholders PowerPoint keeps in sync via
# the Header & Footer dialog's "Date and time" field).

$p = $ppt.ActivePresentation

# ppPlaceholderDate
$ppPlaceholderDate = 16
$oldDate = "14/04/2025"
$newDate = "18/04/2025"

$master = $p.SlideMaster

# Slide master.
foreach ($shp in $master.Shapes) {
    if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate -and $shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# Every slide layout hanging off the (single) slide master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    foreach ($shp in $layout.Shapes) {
        if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Notes master.
foreach ($shp in $p.NotesMaster.Shapes) {
    if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate -and $shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}
